$d = $word.ActiveDocument

$d.Content.Find.Execute("This sample is compatible with the Microsoft Game Development Kit (June 2020)", $true, $false, $false, $false, $false, $true, 1, $false, "This sample is compatible with the Microsoft Game Development Kit (October 2021)", 2)
